# Atualização automática de preços de eletricidade
# Updates row 2 of the Spot_PT sheet with the newly scraped day's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spot_PT")

# Day (date serial)
$ws.Range("A2").Value = 45912

# Hourly prices 0h-1h .. 23h-24h
$ws.Range("B2").Value  = 105
$ws.Range("C2").Value  = 95.03
$ws.Range("D2").Value  = 93.8
$ws.Range("E2").Value  = 93.02
$ws.Range("F2").Value  = 93.02
$ws.Range("G2").Value  = 95.03
$ws.Range("H2").Value  = 105.98
$ws.Range("I2").Value  = 143.15
$ws.Range("J2").Value  = 115.35
$ws.Range("K2").Value  = 93.84
$ws.Range("L2").Value  = 60.97
$ws.Range("M2").Value  = 19.57
$ws.Range("N2").Value  = 5.79
$ws.Range("O2").Value  = 5.76
$ws.Range("P2").Value  = 4.31
$ws.Range("Q2").Value  = 5
$ws.Range("R2").Value  = 4.31
$ws.Range("S2").Value  = 19.68
$ws.Range("T2").Value  = 46.8
$ws.Range("U2").Value  = 93.26000000000001
$ws.Range("V2").Value  = 123.45
$ws.Range("W2").Value  = 143.03
$ws.Range("X2").Value  = 105.01
$ws.Range("Y2").Value  = 97.34999999999999

# Daily average
$ws.Range("Z2").Value = 73.65000000000001

# Slot_4h_max / Slot_4h_price (Slot_4h_max label unchanged: "20h-24h")
$ws.Range("AB2").Value = 117.21

# Slot_2h_frist / Slot_2h_frist_price (Slot_2h_frist label unchanged: "20h-22h")
$ws.Range("AD2").Value = 133.24

# Slot_2h_second / Slot_2h_second_price
$ws.Range("AE2").Value = "6h-8h"
$ws.Range("AF2").Value = 124.56

# Slot_min_price
$ws.Range("AG2").Value = "10h-18h"
